$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell C1 (also updates the table column name automatically)
$ws.Range("C1").Value = "2024-11-15 12:37:41"

# Update data rows: C2, C5, C6 from "Puntual" to "Retardo"
$ws.Range("C2").Value = "Retardo"
$ws.Range("C5").Value = "Retardo"
$ws.Range("C6").Value = "Retardo"
